# TIMECARD-Owen.docx edit:
#  - add a new timecard row for 11/14 ("1" hour, design-model/search-button note)
#    right after the 11/13 row and before the ITERATION III SUBTOTAL row
#  - bump that SUBTOTAL from 24.5 -> 25.5
#  - bump the grand TOTAL from 58.5 -> 59.5

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- locate the "11/13" row; the SUBTOTAL row immediately follows it -----
$subtotalIndex = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    $txt = $row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13)
    if ($txt -eq "11/13") {
        $subtotalIndex = $i + 1
        break
    }
}

# --- insert the new row right before that SUBTOTAL row -------------------
# (re-fetch the row object fresh by index right before use -- holding on to
# a row reference across a mutating call can leave it pointing at the
# wrong row)
$beforeRow = $t.Rows.Item($subtotalIndex)
$newRow = $t.Rows.Add($beforeRow)
$newRow.Cells.Item(1).Range.Text = "11/14"
$newRow.Cells.Item(2).Range.Text = "1"
$newRow.Cells.Item(3).Range.Text = "Revised Design Model and added action to " + [char]8220 + "search" + [char]8221 + " button"

# --- update the ITERATION III SUBTOTAL cell: "24.5" -> "25.5" ------------
# the SUBTOTAL row shifted down by one after the insert
$subtotalRow = $t.Rows.Item($subtotalIndex + 1)
$subtotalCell = $subtotalRow.Cells.Item(2)
$para = $subtotalCell.Range.Paragraphs.Item(1)
$base = $para.Range.Start
# "24.5" => chars: "2"(0) "4"(1) "."(2) "5"(3) ; replace the "4.5" portion with "5.5"
$d.Range($base + 1, $base + 4).Text = "5.5"

# --- update the grand TOTAL paragraph: "TOTAL: 58.5" -> "TOTAL: 59.5" ----
$totalPara = $null
$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith("TOTAL:")) {
        $totalPara = $p
        break
    }
}
$tBase = $totalPara.Range.Start
# "TOTAL: 58.5" => "TOTAL: " is 7 chars, then "5"(7) "8"(8) "."(9) "5"(10)
$d.Range($tBase + 8, $tBase + 11).Text = "9.5"
